# Kandidat_2022_vars.xlsx - malfil edit
# - a5_hovedaktivitet row: C10 snitt_as_num -> snitt_as_num_single, and add a
#   comment in E10 explaining why 2018 comparison is not included.
# - E15:E18 "firedelt skala" comment text tightened (comma instead of ", og").
# - Selection moved to E19 (from E39) and view scrolled back towards the top.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("C10").Value = "snitt_as_num_single"
$ws.Range("E10").Value = "I 2018 fikk ikke kandidatene oppfølgingsspørsmål om grunnen til redusert stilling, og vi har derfor ikke tatt med sammenligning av resultater."

$nyFiredeltTekst = "Før 2022 var svarene gitt med en firedelt skala, vi har derfor ikke tatt med sammenligning av resultater."
$ws.Range("E15").Value = $nyFiredeltTekst
$ws.Range("E16").Value = $nyFiredeltTekst
$ws.Range("E17").Value = $nyFiredeltTekst
$ws.Range("E18").Value = $nyFiredeltTekst

# Scroll the view back up and park the selection on E19 (matches the
# author's saved view state after the edit).
$excel.Goto($ws.Range("B1"), $true)
$ws.Range("E19").Select()
